# Append the new "winners" rows (8-12) to the WinnersData sheet, matching
# the pattern of the existing rows (ID, NAME, THROWS, MONEY columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(6,  "bobi",  4,  20),
    @(7,  "a",     10, 20),
    @(8,  "batko", 2,  20),
    @(9,  "batko", 4,  20),
    @(10, "ivana", 2,  20)
)

$r = 8
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
